$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 14.27711966666667
$ws.Range("H2").Value = 42.831359
$ws.Range("I2").Value = 0.3998945032785361
$ws.Range("J2").Value = 0.399894503278536
$ws.Range("M2").Value = 1.017468666666667
$ws.Range("N2").Value = 3.052406
$ws.Range("O2").Value = 0.2305245503179033
$ws.Range("P2").Value = 0.2305245503179033
$ws.Range("Q2").Value = 14.52652191108378
$ws.Range("R2").Value = 130.738697199754
$ws.Range("S2").Value = 0.09218550054288585
$ws.Range("T2").Value = 0.09218550054288584
$ws.Range("G3").Value = 14.27711966666667
$ws.Range("H3").Value = 42.831359
$ws.Range("I3").Value = 0.3998945032785361
$ws.Range("J3").Value = 0.399894503278536
$ws.Range("O3").Value = 0.5301641128568162
$ws.Range("P3").Value = 0.5301641128568162
$ws.Range("Q3").Value = 33.408331525923
$ws.Range("R3").Value = 300.674983733307
$ws.Range("S3").Value = 0.2120097145669823
$ws.Range("T3").Value = 0.2120097145669822
$ws.Range("G4").Value = 14.27711966666667
$ws.Range("H4").Value = 42.831359
$ws.Range("I4").Value = 0.3998945032785361
$ws.Range("J4").Value = 0.399894503278536
$ws.Range("O4").Value = 0.2393113368252805
$ws.Range("P4").Value = 0.2393113368252805
$ws.Range("Q4").Value = 15.08022192503633
$ws.Range("R4").Value = 135.721997325327
$ws.Range("S4").Value = 0.09569928816866799
$ws.Range("T4").Value = 0.09569928816866798
$ws.Range("I5").Value = 0.3726196691846742
$ws.Range("J5").Value = 0.3726196691846742
$ws.Range("M5").Value = 1.017468666666667
$ws.Range("N5").Value = 3.052406
$ws.Range("O5").Value = 0.2305245503179033
$ws.Range("P5").Value = 0.2305245503179033
$ws.Range("Q5").Value = 13.53573941260644
$ws.Range("R5").Value = 121.821654713458
$ws.Range("S5").Value = 0.08589798167840292
$ws.Range("T5").Value = 0.08589798167840292
$ws.Range("I6").Value = 0.3726196691846742
$ws.Range("J6").Value = 0.3726196691846742
$ws.Range("O6").Value = 0.5301641128568162
$ws.Range("P6").Value = 0.5301641128568162
$ws.Range("S6").Value = 0.1975495763462931
$ws.Range("T6").Value = 0.1975495763462931
$ws.Range("I7").Value = 0.3726196691846742
$ws.Range("J7").Value = 0.3726196691846742
$ws.Range("O7").Value = 0.2393113368252805
$ws.Range("P7").Value = 0.2393113368252805
$ws.Range("S7").Value = 0.08917211115997815
$ws.Range("T7").Value = 0.08917211115997815
$ws.Range("G8").Value = 8.121748
$ws.Range("I8").Value = 0.2274858275367899
$ws.Range("J8").Value = 0.2274858275367898
$ws.Range("M8").Value = 1.017468666666667
$ws.Range("N8").Value = 3.052406
$ws.Range("O8").Value = 0.2305245503179033
$ws.Range("P8").Value = 0.2305245503179033
$ws.Range("Q8").Value = 8.263624108562666
$ws.Range("R8").Value = 74.372616977064
$ws.Range("S8").Value = 0.05244106809661459
$ws.Range("T8").Value = 0.05244106809661459
$ws.Range("G9").Value = 8.121748
$ws.Range("I9").Value = 0.2274858275367899
$ws.Range("J9").Value = 0.2274858275367898
$ws.Range("O9").Value = 0.5301641128568162
$ws.Range("P9").Value = 0.5301641128568162
$ws.Range("S9").Value = 0.1206048219435409
$ws.Range("T9").Value = 0.1206048219435409
$ws.Range("G10").Value = 8.121748
$ws.Range("I10").Value = 0.2274858275367899
$ws.Range("J10").Value = 0.2274858275367898
$ws.Range("O10").Value = 0.2393113368252805
$ws.Range("P10").Value = 0.2393113368252805
$ws.Range("Q10").Value = 8.578604446748001
$ws.Range("R10").Value = 77.20744002073201
$ws.Range("S10").Value = 0.05443993749663439
$ws.Range("T10").Value = 0.05443993749663438
